$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
$ws.Activate()

$values = @(
    @(0.73653905116585106, 0.73170126930228696, 0.74313543198619703, 0.619855248036006),
    @(0.734720869347669,   0.72988308748410602, 0.74113543198619702, 0.61621888439964201),
    @(0.726720869347669,   0.72719336416394798, 0.73446876531952998, 0.61232999551075395),
    @(0.70549864712544696, 0.70386003083061399, 0.70002432087508604, 0.57600736268756103),
    @(0.680215097341897,   0.68179726026784404, 0.69269098754175196, 0.57034069602089399),
    @(0.66947700210380201, 0.67100629036817805, 0.68801252590676898, 0.56906147524167305),
    @(0.58250201203219099, 0.58693483653820999, 0.58142259850845401, 0.46975509395204801),
    @(0.56217381408634604, 0.56007358177989697, 0.54651145318371097, 0.44457320367430803),
    @(0.53421701752831596, 0.53343777477644305, 0.52963928012569705, 0.43007156267639701),
    @(0.50966259624643095, 0.50660730127306797, 0.50385732055235999, 0.40852024687831501),
    @(0.50842283016455903, 0.50548031714608399, 0.50368210152858095, 0.40806047676337198)
)

for ($i = 0; $i -lt 11; $i++) {
    $row = 3 + $i
    $ws.Range("B$row").Value = $values[$i][0]
    $ws.Range("C$row").Value = $values[$i][1]
    $ws.Range("D$row").Value = $values[$i][2]
    $ws.Range("E$row").Value = $values[$i][3]
}

$ws.Range("D15").Select() | Out-Null
